# LTC6811 and SPICE Sim
# Rebuild the BOM table with the updated part list (PCB, DPAK NMOS, SOT23 PMOS,
# PTC thermistor, LDO, 1uF Tantalum capacitor, 13V zener, right-angle crimps /
# housings, 2.5A fuse) and add a grand-Total row (unit total * 7 boards).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- start clean: drop all existing hyperlinks and the old data rows -------
$ws.Hyperlinks.Delete()
$ws.Range("A2:F17").Clear()

# --- header row (unchanged positions, A1 swaps to the "Description" string)-
$ws.Range("A1").Value = "Description"
$ws.Range("B1").Value = "PN"
$ws.Range("C1").Value = "Unit Price"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "Price"
$ws.Range("F1").Value = "URL"

# --- row 2: PCB --------------------------------------------------------
$ws.Range("A2").Value = "PCB"
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = 10
$ws.Range("D2").Value = 1
$ws.Range("E2").Formula = "=C2*D2"

# --- row 3: DPAK NMOS ----------------------------------------------------
$ws.Range("A3").Value = "DPAK NMOS"
$ws.Range("B3").Value = "IRLR110TRPBF"
$ws.Range("C3").Value = 0.65
$ws.Range("D3").Value = 16
$ws.Range("E3").Formula = "=C3*D3"
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.mouser.com/ProductDetail/Vishay-Siliconix/IRLR110TRPBF?qs=sGAEpiMZZMshyDBzk1%2FWi4bBo6KWdDd5qUueznBjLB8%3D ', "", "", "IRLR110TRPBF")

# --- row 4: SOT23 PMOS -----------------------------------------------------
$ws.Range("A4").Value = "SOT23 PMOS"
$ws.Range("B4").Value = "BSS84-7-F"
$ws.Range("C4").Value = 0.066
$ws.Range("D4").Value = 16
$ws.Range("E4").Formula = "=C4*D4"
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.mouser.com/ProductDetail/Nexperia/NX7002AKVL?qs=sGAEpiMZZMshyDBzk1%2FWiwhg%252BJbteykV9w5cubKq8XdkO9kd8b8Ncw%3D%3D ', "", "", "NX7002AKVL")

# --- row 5: PTC thermistor --------------------------------------------------
$ws.Range("A5").Value = "PTC thermistor"
$ws.Range("B5").Value = "PRF15BC102RB6RC"
$ws.Range("C5").Value = 0.187
$ws.Range("D5").Value = 16
$ws.Range("E5").Formula = "=C5*D5"
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.mouser.com/ProductDetail/Murata-Electronics/PRF15BC102RB6RC?qs=sGAEpiMZZMsAYIPNLIqEPNIz%252BwQ02jVW8WQh6KSc6Bc%3D', "", "", "PRF15BC102RB6RC")

# --- row 6: LDO (URL typed as plain text, no live hyperlink) ---------------
$ws.Range("A6").Value = "LDO"
$ws.Range("B6").Value = "LP2951-50DR"
$ws.Range("C6").Value = 0.48
$ws.Range("D6").Value = 2
$ws.Range("E6").Formula = "=C6*D6"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F6").Value = 'https://www.mouser.com/ProductDetail/Texas-Instruments/LP2951-50DR?qs=sGAEpiMZZMsGz1a6aV8DcKyc140gPNQr88A2RaMW6L0%3D'

# --- row 7: 1uF Tantalum capacitor -----------------------------------------
$ws.Range("A7").Value = "1uF Tantalum capacitor"
$ws.Range("B7").Value = "T491A105K035AT"
$ws.Range("C7").Value = 0.235
$ws.Range("D7").Value = 2
$ws.Range("E7").Formula = "=C7*D7"
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.mouser.com/ProductDetail/KEMET/T491A105K035AT?qs=sGAEpiMZZMtZ1n0r9vR22SjDo%2FRaJSx%252BBJL0n6HcLrI%3D ', "", "", "T491A105K035AT")

# --- row 8: 13V zener -------------------------------------------------------
$ws.Range("A8").Value = "13V zener"
$ws.Range("B8").Value = "1SMB5928BT3G"
$ws.Range("C8").Value = 0.316
$ws.Range("D8").Value = 2
$ws.Range("E8").Formula = "=C8*D8"
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.mouser.com/ProductDetail/ON-Semiconductor/1SMB5928BT3G?qs=sGAEpiMZZMtQ8nqTKtFS%2FJ7m6e1KBCguI5yOdiThkFg%3D ', "", "", "1SMB5928BT3G")

# --- row 9: Right angle crimps ----------------------------------------------
$ws.Range("A9").Value = "Right angle crimps"
$ws.Range("B9").Value = "35021-1201"
$ws.Range("C9").Value = 0.03
$ws.Range("D9").Value = 32
$ws.Range("E9").Formula = "=C9*D9"
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.mouser.com/ProductDetail/Molex/35021-1201-Cut-Strip?qs=sGAEpiMZZMs%252BGHln7q6pm%252Bv5BXf4QdrTy6nfkib2RIB4OwsSNmw8Ew%3D%3D ', "", "", "35021-1201")

# --- row 10: Right angle housings -------------------------------------------
$ws.Range("A10").Value = "Right angle housings"
$ws.Range("B10").Value = "35023-0002"
$ws.Range("C10").Value = 0.046
$ws.Range("D10").Value = 16
$ws.Range("E10").Formula = "=C10*D10"
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.mouser.com/ProductDetail/Molex/35023-0002?qs=%2Fha2pyFaduiIAmZSevuPTdYQqsJYR9ufMysb9X1F8SY%3D ', "", "", "35023-0002")

# --- row 11: 2.5A Fuse (URL typed as plain text, no live hyperlink) --------
$ws.Range("A11").Value = "2.5A Fuse"
$ws.Range("B11").Value = "C1F 2.5"
$ws.Range("C11").Value = 0.184
$ws.Range("D11").Value = 16
$ws.Range("E11").Formula = "=C11*D11"
$ws.Range("F11").Style = "Hyperlink"
$ws.Range("F11").Value = 'https://www.mouser.com/ProductDetail/Bel-Fuse/C1F-25?qs=sGAEpiMZZMtxU2g%2F1juGqTi%252BNtKN7qu4Nt8x0xmmQvtWYBQVO0WEHA%3D%3D'

# --- rows 12-14: blank line items (still sum via E column) -----------------
$ws.Range("E12").Formula = "=C12*D12"
$ws.Range("E13").Formula = "=C13*D13"
$ws.Range("E14").Formula = "=C14*D14"

# --- row 18: grand total of the unit column ---------------------------------
$ws.Range("E18").Formula = "=SUM(E2:E14)"

# --- row 19: Total = unit total x 7 boards ----------------------------------
$ws.Range("D19").Value = "Total"
$ws.Range("E19").Formula = "=E18*7"

# --- restore the cursor position seen in the saved workbook ----------------
[void]$ws.Range("A12").Select()
